$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Record Map"
$ws.Range("B7").Value = 99
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("G7").Value = "TEST GRID (TR1)"
$ws.Range("H7").Value = "TEST GRID"

$ws.Range("A8").Select()
